$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 40 for "Destocked - natural land" (non-agricultural land-use)
$ws.Rows("40:40").Insert()

# Insert two new rows after the agricultural-management block (HIR - Beef / HIR - Sheep)
$ws.Rows("47:48").Insert()

# Fill the new row 40
$ws.Range("A40").Value2 = 108
$ws.Range("B40").Value2 = "Destocked - natural land"
$ws.Range("D40").Value2 = "Non-agricultural land-use"
$ws.Range("E40").Value2 = "Non-agricultural land-use"

# Fill the two new rows (47 and 48)
$ws.Range("A47").Value2 = 7
$ws.Range("B47").Value2 = "HIR - Beef"
$ws.Range("E47").Value2 = "Agricultural management"

$ws.Range("A48").Value2 = 8
$ws.Range("B48").Value2 = "HIR - Sheep"
$ws.Range("E48").Value2 = "Agricultural management"

# Restore selection state to match the final author interaction
$ws.Range("B48").Select()
